# Fix hierarchical metrics wrongly computed
# Update the "btop" (column C) and "hitac_filter_qiime" (column H) values
# for rows 2 (precision), 3 (recall) and 4 (f1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.9598784194528877
$ws.Range("C3").Value = 0.9598784194528877
$ws.Range("C4").Value = 0.9598784194528877

$ws.Range("H2").Value = 0.9830278854476348
$ws.Range("H3").Value = 0.9459351801557124
$ws.Range("H4").Value = 0.9641248985015132
